$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Programmatore" (row 6) and "Verificatore" (row 7) hours. ---
# Formulas in column C (=15*B) and the Totale row (8) recalc automatically.
$ws.Range("B6").Value = 20
$ws.Range("B7").Value = 35

# --- Reselect the whole used range, replacing the old multi-area selection. ---
$ws.Range("A1:C8").Select()

# --- Nudge the pie chart a little up/left (pure translation, same size). ---
$co = $ws.ChartObjects().Item(1)
$co.Left = $co.Left - 12
$co.Top = $co.Top - 6.75
